$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.682.53"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "1.847.83"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.030"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.83"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.028"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4384"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3788"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07389"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8818"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.52"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").Value = "1.879.26"
$ws.Range("E12").Value = "  +1.61%  "

$ws.Range("E13").Value = "  +0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.677"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "84.78"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.034"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009063"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.029"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.43"
$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "27.709.44"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.274"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.25"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("D24").Value = "2.107.49"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.049"
$ws.Range("E25").Value = "  +6.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.61"
$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.66"
$ws.Range("E27").Value = "  -0.32%  "

$ws.Range("E28").Value = "  +2.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.312"
$ws.Range("E29").Value = "  +1.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.55"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09060"
$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7732"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.003"
$ws.Range("E34").Value = "  +4.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.548"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("E37").Value = "  +0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01973"
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.839"
$ws.Range("E40").Value = "  +1.68%  "

$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1667"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.850"
$ws.Range("E43").Value = "  +2.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.701"
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.08"
$ws.Range("E45").Value = "  +1.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.69"
$ws.Range("E46").Value = "  +1.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.031"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06568"
$ws.Range("E48").Value = "  +3.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.702"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4690"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.886"
$ws.Range("E51").Value = "  -0.59%  "
